$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "25.574.50"
$ws.Range("E2").Value = "  +1.89%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.666.75"
$ws.Range("E3").Value = "  +0.82%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'236.86"
$ws.Range("E5").Value = "  -0.25%  "

# Row 6 - USDC
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.12%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.2637"
$ws.Range("E8").Value = "  +0.52%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.06155"
$ws.Range("E9").Value = "  +2.87%  "

# Row 10 - TRON
$ws.Range("D10").Value = "'0.07092"
$ws.Range("E10").Value = "  -0.27%  "

# Row 11 - WrappedEther
$ws.Range("D11").Value = "1.667.76"
$ws.Range("E11").Value = "  +0.83%  "

# Row 12 - Solana
$ws.Range("D12").Value = "'14.86"
$ws.Range("E12").Value = "  +2.58%  "

# Row 13 - Polygon
$ws.Range("E13").Value = "  -3.51%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'4.410"
$ws.Range("E14").Value = "  -4.08%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "'74.54"
$ws.Range("E15").Value = "  +1.97%  "

# Row 16 - Dai
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.01%  "

# Row 17 - BinanceUSD
$ws.Range("E17").Value = "  +0.06%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.574.64"
$ws.Range("E18").Value = "  +1.83%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "'0.000006787"
$ws.Range("E19").Value = "  +3.84%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "'11.47"
$ws.Range("E20").Value = "  +0.86%  "

# Row 21 and 22 swap places: Uniswap <-> WrappedliquidstakedEther2.0
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "1.880.85"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'4.466"
$ws.Range("E22").Value = "  +0.91%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "'8.711"
$ws.Range("E23").Value = "  +2.87%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "'5.354"
$ws.Range("E24").Value = "  +1.39%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'134.44"
$ws.Range("E25").Value = "  +1.10%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "'15.11"
$ws.Range("E26").Value = "  +2.63%  "

# Row 27 - Toncoin
$ws.Range("D27").Value = "'1.407"
$ws.Range("E27").Value = "  +1.21%  "

# Row 28 - BitcoinCash
$ws.Range("D28").Value = "'105.02"
$ws.Range("E28").Value = "  +3.14%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'1.695"
$ws.Range("E29").Value = "  -0.68%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "'3.977"
$ws.Range("E30").Value = "  +4.09%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "'3.671"
$ws.Range("E31").Value = "  +4.12%  "

# Row 32 - Stellar
$ws.Range("D32").Value = "'0.07690"
$ws.Range("E32").Value = "  -2.69%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.04358"
$ws.Range("E33").Value = "  -5.15%  "

# Row 34 - Frax
$ws.Range("D34").Value = "'0.9999"
$ws.Range("E34").Value = "  +0.01%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "'2.616"
$ws.Range("E35").Value = "  +0.38%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "'0.6190"
$ws.Range("E36").Value = "  +5.76%  "

# Row 37 - ARBITRUM
$ws.Range("D37").Value = "'0.9525"
$ws.Range("E37").Value = "  +0.78%  "

# Row 38 - MXToken
$ws.Range("D38").Value = "'2.613"
$ws.Range("E38").Value = "  -0.23%  "

# Row 39 - TrustWalletToken
$ws.Range("D39").Value = "'0.8713"
$ws.Range("E39").Value = "  +3.68%  "

# Row 40 - PaxDollar
$ws.Range("D40").Value = "'1.001"
$ws.Range("E40").Value = "  +0.03%  "

# Row 41 - VeChain
$ws.Range("D41").Value = "'0.01516"
$ws.Range("E41").Value = "  -1.70%  "

# Row 42 - RenderToken
$ws.Range("D42").Value = "'1.883"
$ws.Range("E42").Value = "  +2.49%  "

# Row 43 - Quant
$ws.Range("D43").Value = "'97.88"
$ws.Range("E43").Value = "  -0.61%  "

# Row 44 - TheSandbox
$ws.Range("D44").Value = "'0.3780"
$ws.Range("E44").Value = "  +1.80%  "

# Row 45 - FraxShare
$ws.Range("E45").Value = "  -3.09%  "

# Row 46 - Algorand
$ws.Range("D46").Value = "'0.1124"
$ws.Range("E46").Value = "  +0.04%  "

# Row 47 - Aptos
$ws.Range("D47").Value = "'6.232"
$ws.Range("E47").Value = "  +2.81%  "

# Row 48 - Cronos
$ws.Range("D48").Value = "'0.05263"
$ws.Range("E48").Value = "  +2.18%  "

# Row 49 - Elrond
$ws.Range("D49").Value = "'29.54"
$ws.Range("E49").Value = "  +0.82%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "'7.408"
$ws.Range("E50").Value = "  +1.76%  "

# Row 51 - Decentraland -> TrueUSD
$ws.Range("B51").Value = "TrueUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "  +0.21%  "
